$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of B15: remove "localización de los" from the sentence
$ws.Range("B15").Value = "La nutrición en el nivel celular: procesos de endocitosis, fotosíntesis y respiración celular."

# Reflect the active selection on this sheet as B15 (matches the saved sheetView selection)
$ws.Range("B15").Select()
